# Weekly update: insert 3 new price rows (newest week) at the top of the
# data block that starts at row 504, pushing the existing rows down by 3.
# This mirrors the author's workflow of prepending the latest weekly
# observations to the historical series kept in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 504 (existing rows 504:610 shift to 507:613)
$ws.Rows("504:506").Insert()

# New row 504: Pimiento - Zafiro rojo - Primera
$ws.Cells.Item(504, 1).Value = 1
$ws.Cells.Item(504, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(504, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(504, 4).Value = 44694
$ws.Cells.Item(504, 5).Value = 15
$ws.Cells.Item(504, 6).Value = 100112002
$ws.Cells.Item(504, 7).Value = "Pimiento"
$ws.Cells.Item(504, 8).Value = "Zafiro rojo"
$ws.Cells.Item(504, 9).Value = "Primera"
$ws.Cells.Item(504, 10).Value = 120
$ws.Cells.Item(504, 11).Value = 37000
$ws.Cells.Item(504, 12).Value = 38000
$ws.Cells.Item(504, 13).Value = 37500
$ws.Cells.Item(504, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(504, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(504, 16).Value = 2500
$ws.Cells.Item(504, 17).Value = 15
$ws.Cells.Item(504, 18).Value = "Hortaliza"

# New row 505: Pimiento - Zafiro rojo - Segunda
$ws.Cells.Item(505, 1).Value = 1
$ws.Cells.Item(505, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(505, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(505, 4).Value = 44694
$ws.Cells.Item(505, 5).Value = 15
$ws.Cells.Item(505, 6).Value = 100112002
$ws.Cells.Item(505, 7).Value = "Pimiento"
$ws.Cells.Item(505, 8).Value = "Zafiro rojo"
$ws.Cells.Item(505, 9).Value = "Segunda"
$ws.Cells.Item(505, 10).Value = 130
$ws.Cells.Item(505, 11).Value = 34000
$ws.Cells.Item(505, 12).Value = 35000
$ws.Cells.Item(505, 13).Value = 34500
$ws.Cells.Item(505, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(505, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(505, 16).Value = 2300
$ws.Cells.Item(505, 17).Value = 15
$ws.Cells.Item(505, 18).Value = "Hortaliza"

# New row 506: Pimiento - Zafiro rojo - Tercera
$ws.Cells.Item(506, 1).Value = 1
$ws.Cells.Item(506, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(506, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(506, 4).Value = 44694
$ws.Cells.Item(506, 5).Value = 15
$ws.Cells.Item(506, 6).Value = 100112002
$ws.Cells.Item(506, 7).Value = "Pimiento"
$ws.Cells.Item(506, 8).Value = "Zafiro rojo"
$ws.Cells.Item(506, 9).Value = "Tercera"
$ws.Cells.Item(506, 10).Value = 140
$ws.Cells.Item(506, 11).Value = 30000
$ws.Cells.Item(506, 12).Value = 31000
$ws.Cells.Item(506, 13).Value = 30500
$ws.Cells.Item(506, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(506, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(506, 16).Value = 2033
$ws.Cells.Item(506, 17).Value = 15
$ws.Cells.Item(506, 18).Value = "Hortaliza"
